$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "21.3.2018"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = "SYSTEM REQUREMENTS"

$ws.Range("A5").Value = "22.3.218"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = "use case diagram"

$ws.Range("A6").Select()
